$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (old rows 8-12 shift down to 9-13),
# restoring a previously-removed vocabulary entry: 被 / passive / 虛詞
$ws.Rows("8:8").Insert()
$ws.Rows("8:8").RowHeight = 16.5

$ws.Range("A8").Value = "被"
$ws.Range("B8").Value = "passive"
$ws.Range("C8").Value = "虛詞"

# Style the new row with red text and a red box border around A8:C8,
# matching the formatting already used for similar rows in the sheet.
$row8 = $ws.Range("A8:C8")
$row8.Font.Color = 255

$a8 = $ws.Range("A8")
$a8.Borders.Color = 255
$a8.Borders.LineStyle = 1
$a8.Borders.Weight = -4138
$a8.Borders.Item(10).LineStyle = -4142

$b8 = $ws.Range("B8")
$b8.Borders.Color = 255
$b8.Borders.LineStyle = 1
$b8.Borders.Weight = -4138
$b8.Borders.Item(7).LineStyle = -4142
$b8.Borders.Item(10).LineStyle = -4142

$c8 = $ws.Range("C8")
$c8.Borders.Color = 255
$c8.Borders.LineStyle = 1
$c8.Borders.Weight = -4138
$c8.Borders.Item(7).LineStyle = -4142

# Match the selection left behind in the saved workbook
$ws.Range("A8:C8").Select()
